$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
try { Write-Host ("Design.Name=" + $d.Name) } catch { Write-Host ("err1:" + $_.Exception.Message) }
$sm = $d.SlideMaster
try { Write-Host ("sm.Name=" + $sm.Name) } catch { Write-Host ("err2:" + $_.Exception.Message) }
try { $sm.Name = "Office Theme"; Write-Host ("sm.Name now=" + $sm.Name) } catch { Write-Host ("err3:" + $_.Exception.Message) }
